$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.718.95'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '3.779.77'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.80'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.21'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '3.776.15'
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("E13").Value = '  -2.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.06'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '4.413.43'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '3.767.09'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '67.686.43'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.38'
$ws.Range("E18").Value = '  +2.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.02'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.03'
$ws.Range("E21").Value = '  -6.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '457.93'
$ws.Range("E22").Value = '  -1.64%  '
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000153'
$ws.Range("E24").Value = '  +3.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.25'
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.00'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.25'
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.21'
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.68'
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").Value = '3.731.63'
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.100'
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.36'
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.64'
$ws.Range("E44").Value = '  +3.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.15'
$ws.Range("E45").Value = '  +3.05%  '
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.29'
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("E48").Value = '  -1.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '392.06'
$ws.Range("E49").Value = '  +0.81%  '
$ws.Range("E50").Value = '  -4.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.41'
$ws.Range("E51").Value = '  +2.40%  '
